# Add August 2024 ERAS data to both the monthly summary sheet and the
# education-status breakdown sheet, correcting two transposed values
# along the way (G98/G100 on 01_Edu-Status).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 00_Monthly: append row 71 (August 2024 monthly snapshot)
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("00_Monthly")

# Clone formatting from the row above so styles/number formats match.
$wsMonthly.Range("A70:G70").Copy()
$wsMonthly.Range("A71:G71").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsMonthly.Range("A71").Value = 2025
$wsMonthly.Range("B71").Value = 45505
$wsMonthly.Range("C71").Value = 62
$wsMonthly.Range("D71").Value = 766
$wsMonthly.Range("E71").Value = "UK"
$wsMonthly.Range("F71").Value = 12.35
$wsMonthly.Range("G71").Value = "UK"

# ---------------------------------------------------------------------
# 01_Edu-Status: fix the swapped US-DO / IMG mean_apps_program values
# for the 2024-06 period (rows 98 & 100), then append the three new
# August 2024 rows (IMG, US DO, US MD).
# ---------------------------------------------------------------------
$wsEdu = $wb.Worksheets.Item("01_Edu-Status")

$wsEdu.Range("G98").Value = 19.83
$wsEdu.Range("G100").Value = 31

$wsEdu.Range("A110:H110").Copy()
$wsEdu.Range("A113:H113").PasteSpecial(-4122)
$wsEdu.Range("A111:H111").Copy()
$wsEdu.Range("A114:H114").PasteSpecial(-4122)
$wsEdu.Range("A112:H112").Copy()
$wsEdu.Range("A115:H115").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsEdu.Range("A113").Value = 2025
$wsEdu.Range("B113").Value = 45505
$wsEdu.Range("C113").Value = "IMG"
$wsEdu.Range("D113").Value = 50
$wsEdu.Range("E113").Value = 613
$wsEdu.Range("F113").Value = 148
$wsEdu.Range("G113").Value = 12.26
$wsEdu.Range("H113").Value = 4.17

$wsEdu.Range("A114").Value = 2025
$wsEdu.Range("B114").Value = 45505
$wsEdu.Range("C114").Value = "US DO"
$wsEdu.Range("D114").Value = 6
$wsEdu.Range("E114").Value = 20
$wsEdu.Range("F114").Value = 148
$wsEdu.Range("G114").Value = 3.33
$wsEdu.Range("H114").Value = 0.14

$wsEdu.Range("A115").Value = 2025
$wsEdu.Range("B115").Value = 45505
$wsEdu.Range("C115").Value = "US MD"
$wsEdu.Range("D115").Value = 6
$wsEdu.Range("E115").Value = 133
$wsEdu.Range("F115").Value = 148
$wsEdu.Range("G115").Value = 22.17
$wsEdu.Range("H115").Value = 0.9

# ---------------------------------------------------------------------
# View state: 00_Monthly becomes the active tab/selection, and the
# selection on 01_Edu-Status moves to the newly-entered row.
# ---------------------------------------------------------------------
$wsEdu.Activate()
$wsEdu.Range("E115").Select() | Out-Null

$wsMonthly.Activate()
$wsMonthly.Range("F72").Select() | Out-Null
